# Update the PnL-factor calibration estimates with the std-moving-average
# re-fit values (commit: "Aggiornate stime con std moving average").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# AR sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = 0.003255483671362861
$ws.Range("B3").Value = 0.7908842229741908
$ws.Range("B4").Value = 0.1046099214680111
$ws.Range("B5").Value = "[1.0, 0.03564623365737628, -0.04976292463965096, -0.09513492243813464, -0.0780135001702759, 0.04512422601382663, 0.02869501105589732, 0.05104495453266105, 0.010045879271787751, 0.006421806958115386, -0.01390856119209385, -0.05862749443735283, -0.01797437509441243, 0.003360429550374807, 0.061842858814485174, 0.037868044134537114, -0.004120892414323046, -0.029968032245486456, -0.046788926839067546, -0.004517349322110577]"

# ---------------------------------------------------------------------
# SETAR sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B4").Value = -0.1701720019737109
$ws.Range("B5").Value = 0.6038208696320441
$ws.Range("B6").Value = 0.06207413713487229
$ws.Range("B7").Value = 0.1736986288907789
$ws.Range("B8").Value = 0.5888753153466679
$ws.Range("B9").Value = 0.06453489138270903
$ws.Range("B10").Value = "[1.0, 0.1005583244946701, -0.0033006177235156866, 0.012032751791994898, 0.0038070452406806727, 0.0034889842893371997, -0.012894027705417617, 0.024529143340206748, -0.006286581110110385, -0.029482700629761497, -0.007907859454282897, -0.03164763329861033, 0.0035084880707617304, -0.007644827993641003, 0.01286451460552699, 0.010663892768255044, 0.008296078400547411, -0.005567066780091306, -0.01724534853220375, 0.009857257540119414]"

# ---------------------------------------------------------------------
# GARCH sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = -0.00005104153231234426
$ws.Range("B3").Value = 0.111236535563846
$ws.Range("B4").Value = 0.04715544375455139
$ws.Range("B5").Value = 0.00000003120092914657768
$ws.Range("B6").Value = "[1.0, -0.0012218618642324186, -0.03854571935315319, -0.08050302102553594, -0.08948667694000302, 0.05789108169828436, 0.02015598037768332, 0.050920649085659235, 0.007879646456085522, 0.012774177981537573, -0.027219973500138175, -0.05823674657093854, -0.013438638779686315, 0.003469158622862317, 0.061403784376428884, 0.03216206813663422, 0.009551467403641708, -0.03982963560937029, -0.04133427500069295, -0.004183119619358411]"

# ---------------------------------------------------------------------
# TARCH sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = 0.0002355363274201636
$ws.Range("B3").Value = 0.1112582581923891
$ws.Range("B4").Value = 0.05408906711830145
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = "[0.9999999999999998, -0.0010073791534656582, -0.03853290282442895, -0.08060476262221215, -0.0895434994307278, 0.057636838424661524, 0.02052564447657671, 0.05109678197350378, 0.007936922665196218, 0.012448539965017054, -0.02742792392836405, -0.058346202226263306, -0.013014366150564545, 0.0036206434956232406, 0.06128489758827131, 0.032106043404878395, 0.00967159383880573, -0.039888243343761114, -0.041443412775109795, -0.004050248934138563]"
$ws.Range("B7").Value = -0.01431547116041601

# ---------------------------------------------------------------------
# AR_TARCH sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = 0.003179021498198714
$ws.Range("B3").Value = 0.1006705903290248
$ws.Range("B4").Value = 0.04796953377191113
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = "[1.0, 0.00039690472789153666, -0.05006247197088025, -0.09357214280955617, -0.0781956366378188, 0.04671096906923515, 0.026506634098402388, 0.05124155187022357, 0.007988591435446864, 0.006588003242049092, -0.014654973911823566, -0.058753404853894714, -0.016109912803037252, 0.0023788615285235157, 0.06171479993223185, 0.03669935493715551, -0.0038563817503224577, -0.029714641819497543, -0.04707105663136452, -0.004073870981220352]"
$ws.Range("B7").Value = -0.02130791236529418
$ws.Range("B9").Value = 0.7905020366613943
